# "delete first new line" - remove the "New item" paragraph entirely
# (its text plus its paragraph mark). Word keeps the hidden "_GoBack"
# bookmark pinned to the location of the most recent edit, so after the
# deletion we re-add it (by name, which moves an existing bookmark)
# collapsed at the start of the paragraph that now follows the deletion
# point ("Second new item").

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd("`r", "`n", "`x07")
    if ($t -eq "New item") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $deletePos = $target.Range.Start
    $target.Range.Delete()

    $goBackRange = $d.Range($deletePos, $deletePos)
    $d.Bookmarks.Add("_GoBack", $goBackRange)
}
